# Auto-generated Excel COM-interop script applying cell value changes
# derived from the canonical OOXML diff of Sheets/Balmung_Profits.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I4").Value = 554.5
$ws.Range("K4").Value = 554.5
$ws.Range("M4").Value = -440.5
$ws.Range("H4").Value = 1012.7143
$ws.Range("I15").Value = 950.8298
$ws.Range("M15").Value = -2683.4894
$ws.Range("H15").Value = 950.8298
$ws.Range("K15").Value = 2852.4894
$ws.Range("J18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("L18").Value = 0
$ws.Range("H18").Value = 3194.25
$ws.Range("K31").Value = 2251.5
$ws.Range("H31").Value = 750.5
$ws.Range("I31").Value = 750.5
$ws.Range("M31").Value = -2021.5
$ws.Range("H38").Value = 4322.4116
$ws.Range("J38").Value = 8999.25
$ws.Range("N38").Value = -27741.75
$ws.Range("L38").Value = 26997.75
$ws.Range("K40").Value = 4735.28
$ws.Range("I40").Value = 4735.28
$ws.Range("H40").Value = 4659.697
$ws.Range("M40").Value = -4560.28
$ws.Range("L58").Value = 6051
$ws.Range("N58").Value = -6351
$ws.Range("J58").Value = 2017
$ws.Range("H58").Value = 483.64706
$ws.Range("H62").Value = 6924
$ws.Range("K62").Value = 5265
$ws.Range("M62").Value = -4641
$ws.Range("I62").Value = 5265
$ws.Range("H65").Value = 6924
$ws.Range("M65").Value = -23205
$ws.Range("K65").Value = 26325
$ws.Range("I65").Value = 5265
$ws.Range("N69").Value = -141329.75
$ws.Range("J69").Value = 46527.25
$ws.Range("H69").Value = 42221.6
$ws.Range("L69").Value = 139581.75
$ws.Range("L72").Value = 418745.25
$ws.Range("N72").Value = -427481.25
$ws.Range("J72").Value = 46527.25
$ws.Range("H72").Value = 42221.6
$ws.Range("L75").Value = 46303.5
$ws.Range("J75").Value = 46303.5
$ws.Range("H75").Value = 46303.5
$ws.Range("N75").Value = -48175.5
$ws.Range("J78").Value = 46303.5
$ws.Range("L78").Value = 138910.5
$ws.Range("N78").Value = -148270.5
$ws.Range("H78").Value = 46303.5
$ws.Range("N86").Value = -83337736
$ws.Range("J86").Value = 83335490
$ws.Range("I86").Value = 76926296
$ws.Range("H86").Value = 78434344
$ws.Range("M86").Value = -76925173
$ws.Range("L86").Value = 83335490
$ws.Range("K86").Value = 76926296
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 100000
$ws.Range("N87").Value = -102496
$ws.Range("L87").Value = 100000
$ws.Range("J89").Value = 83335490
$ws.Range("M89").Value = -384625864
$ws.Range("L89").Value = 416677450
$ws.Range("N89").Value = -416688682
$ws.Range("H89").Value = 78434344
$ws.Range("K89").Value = 384631480
$ws.Range("I89").Value = 76926296
$ws.Range("N90").Value = -312480
$ws.Range("H90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("J90").Value = 100000
$ws.Range("H98").Value = 1881.9615
$ws.Range("M98").Value = 447.55
$ws.Range("K98").Value = 1050.45
$ws.Range("I98").Value = 1050.45
$ws.Range("H111").Value = 4339.0835
$ws.Range("K111").Value = 12621.3
$ws.Range("M111").Value = -9554.300000000001
$ws.Range("I111").Value = 4207.1
$ws.Range("N113").Value = -16281.667
$ws.Range("L113").Value = 9773.666999999999
$ws.Range("J113").Value = 9773.666999999999
$ws.Range("H113").Value = 7761.1665
$ws.Range("M122").Value = -701.3500000000004
$ws.Range("H122").Value = 1881.9615
$ws.Range("K122").Value = 3151.35
$ws.Range("I122").Value = 1050.45
$ws.Range("K132").Value = 159575.1
$ws.Range("I132").Value = 53191.7
$ws.Range("H132").Value = 41528.152
$ws.Range("M132").Value = -157045.1
$ws.Range("J135").Value = 2197.1538
$ws.Range("L135").Value = 19774.3842
$ws.Range("K135").Value = 19724.04
$ws.Range("M135").Value = -17189.04
$ws.Range("H135").Value = 2193.4736
$ws.Range("I135").Value = 2191.56
$ws.Range("N135").Value = -24844.3842
$ws.Range("K138").Value = 39664.236
$ws.Range("I138").Value = 13221.412
$ws.Range("L138").Value = 17581.5792
$ws.Range("M138").Value = -34524.236
$ws.Range("J138").Value = 5860.5264
$ws.Range("N138").Value = -27861.5792
$ws.Range("H138").Value = 8135.709

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H6").Value = 0
$ws.Range("J32").Value = 29591.076
$ws.Range("K32").Value = 9286708
$ws.Range("N32").Value = -30165.076
$ws.Range("H32").Value = 6879858
$ws.Range("I32").Value = 9286708
$ws.Range("L32").Value = 29591.076
$ws.Range("M32").Value = -9286421
$ws.Range("M34").Value = -212395.67
$ws.Range("I34").Value = 212666.67
$ws.Range("H34").Value = 212666.67
$ws.Range("K34").Value = 212666.67
$ws.Range("H45").Value = 2823.5264
$ws.Range("I45").Value = 2975.5454
$ws.Range("M45").Value = -2598.5454
$ws.Range("K45").Value = 2975.5454
$ws.Range("I61").Value = 12254.3125
$ws.Range("H61").Value = 1938616.2
$ws.Range("M61").Value = -12042.3125
$ws.Range("K61").Value = 12254.3125
$ws.Range("M74").Value = -3007.4324
$ws.Range("H74").Value = 432044.2
$ws.Range("I74").Value = 3881.4324
$ws.Range("K74").Value = 3881.4324
$ws.Range("M77").Value = -15039.162
$ws.Range("K77").Value = 19407.162
$ws.Range("I77").Value = 3881.4324
$ws.Range("H77").Value = 432044.2
$ws.Range("J88").Value = 3085.6667
$ws.Range("N88").Value = -3897.6667
$ws.Range("L88").Value = 3085.6667
$ws.Range("H88").Value = 3040.2727
$ws.Range("N91").Value = -5893.6667
$ws.Range("J91").Value = 3085.6667
$ws.Range("L91").Value = 3085.6667
$ws.Range("H91").Value = 3040.2727
$ws.Range("I97").Value = 4972.2085
$ws.Range("H97").Value = 4132.485
$ws.Range("M97").Value = -4476.2085
$ws.Range("K97").Value = 4972.2085
$ws.Range("K102").Value = 3941.353
$ws.Range("M102").Value = -2319.353
$ws.Range("H102").Value = 4319.227
$ws.Range("I102").Value = 3941.353
$ws.Range("M122").Value = -1951.3333
$ws.Range("H122").Value = 1468.5454
$ws.Range("K122").Value = 4401.3333
$ws.Range("I122").Value = 1467.1111
$ws.Range("K132").Value = 13432.0431
$ws.Range("I132").Value = 4477.3477
$ws.Range("H132").Value = 4743.086
$ws.Range("M132").Value = -10902.0431
$ws.Range("L134").Value = 127983.8
$ws.Range("J134").Value = 127983.8
$ws.Range("N134").Value = -138123.8
$ws.Range("H134").Value = 127983.8
$ws.Range("H136").Value = 1938616.2
$ws.Range("M136").Value = -34212.9375
$ws.Range("K136").Value = 36762.9375
$ws.Range("I136").Value = 12254.3125
$ws.Range("J140").Value = 183611.55
$ws.Range("H140").Value = 183611.55
$ws.Range("N140").Value = -193971.55
$ws.Range("L140").Value = 183611.55

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M22").Value = -1212.625
$ws.Range("I22").Value = 1385.625
$ws.Range("K22").Value = 1385.625
$ws.Range("H22").Value = 1385.625
$ws.Range("I86").Value = 2241.1667
$ws.Range("H86").Value = 3287.0833
$ws.Range("M86").Value = -1118.1667
$ws.Range("K86").Value = 2241.1667
$ws.Range("M89").Value = -5589.833500000001
$ws.Range("H89").Value = 3287.0833
$ws.Range("K89").Value = 11205.8335
$ws.Range("I89").Value = 2241.1667
$ws.Range("K94").Value = 2950.8125
$ws.Range("H94").Value = 3419.524
$ws.Range("M94").Value = -2499.8125
$ws.Range("I94").Value = 2950.8125
$ws.Range("M96").Value = -16273.4
$ws.Range("H96").Value = 19019.4
$ws.Range("I96").Value = 19019.4
$ws.Range("K96").Value = 19019.4
$ws.Range("M99").Value = -9431.538
$ws.Range("L99").Value = 480
$ws.Range("N99").Value = -3476
$ws.Range("J99").Value = 480
$ws.Range("H99").Value = 9536.267
$ws.Range("K99").Value = 10929.538
$ws.Range("I99").Value = 10929.538
$ws.Range("M107").Value = -10803.765
$ws.Range("I107").Value = 12723.765
$ws.Range("H107").Value = 10831.681
$ws.Range("K107").Value = 12723.765
$ws.Range("M134").Value = -111265.899
$ws.Range("K134").Value = 113800.899
$ws.Range("H134").Value = 20481780
$ws.Range("I134").Value = 37933.633
$ws.Range("J139").Value = 32000
$ws.Range("H139").Value = 32000
$ws.Range("N139").Value = -42280
$ws.Range("L139").Value = 32000

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M2").Value = -1150.3334
$ws.Range("I2").Value = 1263.3334
$ws.Range("H2").Value = 1263.3334
$ws.Range("K2").Value = 1263.3334
$ws.Range("H16").Value = 2675.6
$ws.Range("L16").Value = 2549.4285
$ws.Range("N16").Value = -3123.4285
$ws.Range("J16").Value = 2549.4285
$ws.Range("M22").Value = -204.381
$ws.Range("L22").Value = 13566.333
$ws.Range("I22").Value = 554.381
$ws.Range("J22").Value = 13566.333
$ws.Range("N22").Value = -14266.333
$ws.Range("K22").Value = 554.381
$ws.Range("H22").Value = 2180.875
$ws.Range("J31").Value = 6231.36
$ws.Range("K31").Value = 3583.4614
$ws.Range("H31").Value = 5325.5
$ws.Range("L31").Value = 6231.36
$ws.Range("I31").Value = 3583.4614
$ws.Range("M31").Value = -3288.4614
$ws.Range("N31").Value = -6821.36
$ws.Range("M34").Value = -3381.4614
$ws.Range("J34").Value = 6231.36
$ws.Range("I34").Value = 3583.4614
$ws.Range("L34").Value = 6231.36
$ws.Range("N34").Value = -6635.36
$ws.Range("H34").Value = 5325.5
$ws.Range("K34").Value = 3583.4614
$ws.Range("I41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("M41").ClearContents()
$ws.Range("K41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H51").Value = 15000
$ws.Range("L58").Value = 4681.5
$ws.Range("M58").Value = -2258.375
$ws.Range("K58").Value = 2461.375
$ws.Range("N58").Value = -5087.5
$ws.Range("I58").Value = 2461.375
$ws.Range("J58").Value = 4681.5
$ws.Range("H58").Value = 3412.8572
$ws.Range("K60").Value = 7667
$ws.Range("I60").Value = 7667
$ws.Range("H60").Value = 7667
$ws.Range("M60").Value = -7156
$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("L61").Value = 0
$ws.Range("N86").Value = -17405.909
$ws.Range("J86").Value = 15159.909
$ws.Range("I86").Value = 5349.467
$ws.Range("H86").Value = 9500.038
$ws.Range("M86").Value = -4226.467
$ws.Range("L86").Value = 15159.909
$ws.Range("K86").Value = 5349.467
$ws.Range("J89").Value = 15159.909
$ws.Range("M89").Value = -21131.335
$ws.Range("L89").Value = 75799.545
$ws.Range("N89").Value = -87031.545
$ws.Range("H89").Value = 9500.038
$ws.Range("K89").Value = 26747.335
$ws.Range("I89").Value = 5349.467
$ws.Range("J96").Value = 34633
$ws.Range("N96").Value = -40125
$ws.Range("H96").Value = 34633
$ws.Range("L96").Value = 34633
$ws.Range("K105").Value = 1409
$ws.Range("M105").Value = 338
$ws.Range("I105").Value = 1409
$ws.Range("H105").Value = 1938.5555
$ws.Range("M107").Value = -2656.9
$ws.Range("I107").Value = 4576.9
$ws.Range("H107").Value = 5588.2666
$ws.Range("K107").Value = 4576.9
$ws.Range("N113").Value = -6889.4285
$ws.Range("L113").Value = 2549.4285
$ws.Range("J113").Value = 2549.4285
$ws.Range("H113").Value = 2675.6
$ws.Range("M134").Value = -4283.4786
$ws.Range("K134").Value = 6818.4786
$ws.Range("H134").Value = 2932.818
$ws.Range("I134").Value = 2272.8262
$ws.Range("J136").Value = 4681.5
$ws.Range("L136").Value = 14044.5
$ws.Range("H136").Value = 3412.8572
$ws.Range("M136").Value = -4834.125
$ws.Range("K136").Value = 7384.125
$ws.Range("I136").Value = 2461.375
$ws.Range("N136").Value = -19144.5
$ws.Range("L141").Value = 511412.78
$ws.Range("H141").Value = 529782.3
$ws.Range("N141").Value = -521772.78
$ws.Range("J141").Value = 511412.78

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K3").Value = 1767
$ws.Range("I3").Value = 589
$ws.Range("H3").Value = 7538.3335
$ws.Range("M3").Value = -1655
$ws.Range("I4").Value = 3536456.8
$ws.Range("K4").Value = 10609370.4
$ws.Range("M4").Value = -10609258.4
$ws.Range("H4").Value = 2140626.2
$ws.Range("M17").Value = -237.5
$ws.Range("H17").Value = 368.6
$ws.Range("K17").Value = 406.5
$ws.Range("I17").Value = 135.5
$ws.Range("L25").Value = 8994.75
$ws.Range("H25").Value = 2728.5
$ws.Range("N25").Value = -9332.75
$ws.Range("J25").Value = 2998.25
$ws.Range("J30").Value = 2998.25
$ws.Range("H30").Value = 2728.5
$ws.Range("N30").Value = -9198.75
$ws.Range("L30").Value = 8994.75
$ws.Range("J34").Value = 472.91306
$ws.Range("L34").Value = 1418.73918
$ws.Range("N34").Value = -1586.73918
$ws.Range("H34").Value = 469.83334
$ws.Range("J88").Value = 10119.435
$ws.Range("N88").Value = -31214.305
$ws.Range("I88").Value = 3997.5
$ws.Range("L88").Value = 30358.305
$ws.Range("H88").Value = 9629.68
$ws.Range("M88").Value = -11564.5
$ws.Range("K88").Value = 11992.5
$ws.Range("N91").Value = -33322.305
$ws.Range("M91").Value = -10510.5
$ws.Range("K91").Value = 11992.5
$ws.Range("J91").Value = 10119.435
$ws.Range("I91").Value = 3997.5
$ws.Range("L91").Value = 30358.305
$ws.Range("H91").Value = 9629.68
$ws.Range("H92").Value = 166.13333
$ws.Range("L92").Value = 386.33331
$ws.Range("J92").Value = 128.77777
$ws.Range("N92").Value = -2882.33331
$ws.Range("I93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -22344
$ws.Range("J93").Value = 6200
$ws.Range("H93").Value = 6200
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 18600
$ws.Range("H137").Value = 6206.72
$ws.Range("N137").Value = -29620.3842
$ws.Range("L137").Value = 19420.3842
$ws.Range("J137").Value = 6473.4614
$ws.Range("L138").Value = 13000.0005
$ws.Range("J138").Value = 4333.3335
$ws.Range("N138").Value = -23280.0005
$ws.Range("H138").Value = 4130.357
$ws.Range("K139").Value = 19741375.5
$ws.Range("M139").Value = -19736235.5
$ws.Range("I139").Value = 6580458.5
$ws.Range("H139").Value = 4169509.8
$ws.Range("I141").Value = 3504.9167
$ws.Range("K141").Value = 10514.7501
$ws.Range("M141").Value = -5334.750100000001
$ws.Range("H141").Value = 3504.9167

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M21").Value = -4525.25
$ws.Range("H21").Value = 4698.25
$ws.Range("I21").Value = 4698.25
$ws.Range("K21").Value = 4698.25
$ws.Range("K30").Value = 4698.25
$ws.Range("I30").Value = 4698.25
$ws.Range("H30").Value = 4698.25
$ws.Range("M30").Value = -4593.25
$ws.Range("I80").Value = 184837.14
$ws.Range("H80").Value = 9112179
$ws.Range("M80").Value = -183839.14
$ws.Range("K80").Value = 184837.14
$ws.Range("I83").Value = 184837.14
$ws.Range("K83").Value = 924185.7000000001
$ws.Range("M83").Value = -919193.7000000001
$ws.Range("H83").Value = 9112179
$ws.Range("H87").Value = 100354
$ws.Range("J87").Value = 100354
$ws.Range("N87").Value = -102850
$ws.Range("L87").Value = 100354
$ws.Range("J88").Value = 102684.664
$ws.Range("N88").Value = -103586.664
$ws.Range("L88").Value = 102684.664
$ws.Range("H88").Value = 102684.664
$ws.Range("N90").Value = -313542
$ws.Range("H90").Value = 100354
$ws.Range("L90").Value = 301062
$ws.Range("J90").Value = 100354
$ws.Range("N91").Value = -105804.664
$ws.Range("J91").Value = 102684.664
$ws.Range("L91").Value = 102684.664
$ws.Range("H91").Value = 102684.664
$ws.Range("I97").Value = 939.2917
$ws.Range("H97").Value = 4123.6665
$ws.Range("M97").Value = -443.2917
$ws.Range("K97").Value = 939.2917
$ws.Range("J107").Value = 220
$ws.Range("L107").Value = 220
$ws.Range("N107").Value = -4060
$ws.Range("M107").Value = -65109.266
$ws.Range("I107").Value = 67029.266
$ws.Range("H107").Value = 47940.906
$ws.Range("K107").Value = 67029.266
$ws.Range("H109").Value = 97227.5
$ws.Range("J109").Value = 97227.5
$ws.Range("L109").Value = 97227.5
$ws.Range("N109").Value = -99307.5
$ws.Range("M113").Value = -2830
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("H113").Value = 5000
$ws.Range("M122").Value = -14370.7861
$ws.Range("H122").Value = 4562.421
$ws.Range("K122").Value = 16820.7861
$ws.Range("I122").Value = 5606.9287
$ws.Range("K132").Value = 9125.694
$ws.Range("L132").Value = 37600035
$ws.Range("I132").Value = 3041.898
$ws.Range("J132").Value = 12533345
$ws.Range("H132").Value = 3885671
$ws.Range("M132").Value = -6595.694
$ws.Range("N132").Value = -37605095
$ws.Range("J136").Value = 45087.09
$ws.Range("L136").Value = 135261.27
$ws.Range("H136").Value = 45087.09
$ws.Range("N136").Value = -140361.27

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K9").Value = 2187
$ws.Range("I9").Value = 2187
$ws.Range("M9").Value = -1963
$ws.Range("H9").Value = 2589.6
$ws.Range("H16").Value = 984.25
$ws.Range("L16").Value = 1199
$ws.Range("K16").Value = 934.6923
$ws.Range("N16").Value = -1539
$ws.Range("I16").Value = 934.6923
$ws.Range("M16").Value = -764.6923
$ws.Range("J16").Value = 1199
$ws.Range("M22").Value = -1321.1818
$ws.Range("L22").Value = 4624.4443
$ws.Range("I22").Value = 1616.1818
$ws.Range("J22").Value = 4624.4443
$ws.Range("N22").Value = -5214.4443
$ws.Range("K22").Value = 1616.1818
$ws.Range("H22").Value = 3483.3794
$ws.Range("N27").Value = -4838.4443
$ws.Range("M27").Value = -1509.1818
$ws.Range("J27").Value = 4624.4443
$ws.Range("L27").Value = 4624.4443
$ws.Range("K27").Value = 1616.1818
$ws.Range("I27").Value = 1616.1818
$ws.Range("H27").Value = 3483.3794
$ws.Range("N46").Value = -5774
$ws.Range("H46").Value = 13052.223
$ws.Range("M46").Value = -15808.154
$ws.Range("J46").Value = 5398
$ws.Range("L46").Value = 5398
$ws.Range("I46").Value = 15996.154
$ws.Range("K46").Value = 15996.154
$ws.Range("I61").Value = 1475.3636
$ws.Range("H61").Value = 2368.6296
$ws.Range("M61").Value = -1273.3636
$ws.Range("K61").Value = 1475.3636
$ws.Range("N64").Value = -91334
$ws.Range("H64").Value = 90884
$ws.Range("J64").Value = 90884
$ws.Range("L64").Value = 90884
$ws.Range("N67").Value = -92444
$ws.Range("H67").Value = 90884
$ws.Range("J67").Value = 90884
$ws.Range("L67").Value = 90884
$ws.Range("M113").Value = 694.6364000000001
$ws.Range("I113").Value = 1475.3636
$ws.Range("K113").Value = 1475.3636
$ws.Range("H113").Value = 2368.6296
$ws.Range("M122").Value = -12434.7139
$ws.Range("L122").Value = 17456.7
$ws.Range("J122").Value = 5818.9
$ws.Range("N122").Value = -22356.7
$ws.Range("H122").Value = 5318.7915
$ws.Range("K122").Value = 14884.7139
$ws.Range("I122").Value = 4961.5713
$ws.Range("K132").Value = 16654.125
$ws.Range("I132").Value = 5551.375
$ws.Range("H132").Value = 12440.55
$ws.Range("M132").Value = -14124.125
$ws.Range("L133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("H133").Value = 80000
$ws.Range("N133").Value = -85060
$ws.Range("H136").Value = 7214.8096
$ws.Range("M136").Value = -20321.769
$ws.Range("K136").Value = 22871.769
$ws.Range("I136").Value = 7623.923
$ws.Range("L138").Value = 176499.25
$ws.Range("J138").Value = 176499.25
$ws.Range("N138").Value = -186779.25
$ws.Range("H138").Value = 176499.25
$ws.Range("J140").Value = 122397.6
$ws.Range("H140").Value = 122397.6
$ws.Range("N140").Value = -132757.6
$ws.Range("L140").Value = 122397.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 60420
$ws.Range("L16").Value = 60420
$ws.Range("N16").Value = -61004
$ws.Range("J16").Value = 60420
$ws.Range("H62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("J63").Value = 40416.332
$ws.Range("N63").Value = -41664.332
$ws.Range("L63").Value = 40416.332
$ws.Range("H63").Value = 38937.25
$ws.Range("N64").Value = -40496
$ws.Range("H64").Value = 37499.75
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("J65").Value = 0
$ws.Range("N66").Value = -127488.996
$ws.Range("H66").Value = 38937.25
$ws.Range("L66").Value = 121248.996
$ws.Range("J66").Value = 40416.332
$ws.Range("N67").Value = -41716
$ws.Range("L67").Value = 40000
$ws.Range("H67").Value = 37499.75
$ws.Range("J67").Value = 40000
$ws.Range("M81").Value = -11341887
$ws.Range("I81").Value = 5671474
$ws.Range("H81").Value = 5671474
$ws.Range("K81").Value = 11342948
$ws.Range("H84").Value = 5671474
$ws.Range("I84").Value = 5671474
$ws.Range("M84").Value = -56709436
$ws.Range("K84").Value = 56714740
$ws.Range("H92").Value = 59713.07
$ws.Range("L92").Value = 59713.07
$ws.Range("J92").Value = 59713.07
$ws.Range("N92").Value = -64705.07
$ws.Range("M122").Value = -4358.2225
$ws.Range("H122").Value = 2219.0667
$ws.Range("K122").Value = 6808.2225
$ws.Range("I122").Value = 2269.4075
$ws.Range("L132").Value = 8385.882599999999
$ws.Range("J132").Value = 2795.2942
$ws.Range("H132").Value = 49952.332
$ws.Range("N132").Value = -13445.8826
$ws.Range("J136").Value = 13151
$ws.Range("L136").Value = 39453
$ws.Range("H136").Value = 9620.75
$ws.Range("N136").Value = -44553
$ws.Range("J140").Value = 121395.4
$ws.Range("H140").Value = 121395.4
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -131755.4
$ws.Range("L140").Value = 121395.4

Write-Host "Applied all cell updates."